# Update the team-specific transition-probability matrix on Sheet1 with
# recomputed values reflecting additional simulated games / refreshed stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1964285714285714
$ws.Range("C2").Value = 0.5446428571428571
$ws.Range("J2").Value = 0.01785714285714286
$ws.Range("P2").Value = 0.1383928571428572
$ws.Range("S2").Value = 0.1026785714285714
$ws.Range("B3").Value = 0.01550387596899225
$ws.Range("C3").Value = 0.007751937984496124
$ws.Range("J3").Value = 0.05426356589147287
$ws.Range("P3").Value = 0.7286821705426356
$ws.Range("S3").Value = 0.1937984496124031
$ws.Range("J4").Value = 0.0303030303030303
$ws.Range("P4").Value = 0.696969696969697
$ws.Range("S4").Value = 0.2727272727272727
$ws.Range("B6").Value = 0.05555555555555555
$ws.Range("D6").Value = 0.004629629629629629
$ws.Range("F6").Value = 0.05092592592592592
$ws.Range("J6").Value = 0.25
$ws.Range("O6").Value = 0.01388888888888889
$ws.Range("Q6").Value = 0.1574074074074074
$ws.Range("R6").Value = 0.1018518518518518
$ws.Range("S6").Value = 0.3657407407407408
$ws.Range("B7").Value = 0.06040268456375839
$ws.Range("D7").Value = 0.01342281879194631
$ws.Range("E7").Value = 0.006711409395973154
$ws.Range("F7").Value = 0.04697986577181208
$ws.Range("J7").Value = 0.2080536912751678
$ws.Range("O7").Value = 0.02013422818791946
$ws.Range("Q7").Value = 0.1006711409395973
$ws.Range("R7").Value = 0.1073825503355705
$ws.Range("S7").Value = 0.436241610738255
$ws.Range("B8").Value = 0.06091370558375635
$ws.Range("D8").Value = 0.005076142131979695
$ws.Range("F8").Value = 0.04822335025380711
$ws.Range("J8").Value = 0.09644670050761421
$ws.Range("O8").Value = 0.007614213197969543
$ws.Range("Q8").Value = 0.2106598984771574
$ws.Range("R8").Value = 0.1065989847715736
$ws.Range("S8").Value = 0.4644670050761421
$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("D9").Value = 0.01136363636363636
$ws.Range("F9").Value = 0.07954545454545454
$ws.Range("J9").Value = 0.1363636363636364
$ws.Range("O9").Value = 0.01136363636363636
$ws.Range("Q9").Value = 0.1363636363636364
$ws.Range("R9").Value = 0.1098484848484848
$ws.Range("S9").Value = 0.4242424242424243
$ws.Range("B10").Value = 0.08073115003808073
$ws.Range("D10").Value = 0.0198019801980198
$ws.Range("E10").Value = 0.0007616146230007616
$ws.Range("F10").Value = 0.08149276466108149
$ws.Range("J10").Value = 0.1264280274181264
$ws.Range("O10").Value = 0.01066260472201066
$ws.Range("Q10").Value = 0.1957349581111957
$ws.Range("R10").Value = 0.09748667174409749
$ws.Range("S10").Value = 0.3869002284843869
$ws.Range("G11").Value = 0.1192660550458716
$ws.Range("J11").Value = 0.07339449541284404
$ws.Range("K11").Value = 0.1697247706422018
$ws.Range("L11").Value = 0.6146788990825688
$ws.Range("S11").Value = 0.02293577981651376
$ws.Range("G12").Value = 0.7163120567375887
$ws.Range("J12").Value = 0.1843971631205674
$ws.Range("K12").Value = 0.02836879432624113
$ws.Range("L12").Value = 0.05673758865248227
$ws.Range("S12").Value = 0.01418439716312057
$ws.Range("G13").Value = 0.8620689655172413
$ws.Range("J13").Value = 0.103448275862069
$ws.Range("S13").Value = 0.03448275862068965
$ws.Range("F15").Value = 0.02173913043478261
$ws.Range("H15").Value = 0.1576086956521739
$ws.Range("I15").Value = 0.1141304347826087
$ws.Range("J15").Value = 0.4293478260869565
$ws.Range("K15").Value = 0.04891304347826087
$ws.Range("M15").Value = 0.01630434782608696
$ws.Range("O15").Value = 0.04891304347826087
$ws.Range("S15").Value = 0.1630434782608696
$ws.Range("F16").Value = 0.00684931506849315
$ws.Range("H16").Value = 0.136986301369863
$ws.Range("I16").Value = 0.0821917808219178
$ws.Range("J16").Value = 0.4726027397260274
$ws.Range("K16").Value = 0.07534246575342465
$ws.Range("M16").Value = 0.0136986301369863
$ws.Range("N16").Value = 0.00684931506849315
$ws.Range("O16").Value = 0.07534246575342465
$ws.Range("S16").Value = 0.1301369863013699
$ws.Range("F17").Value = 0.01187648456057007
$ws.Range("H17").Value = 0.1401425178147268
$ws.Range("I17").Value = 0.1187648456057007
$ws.Range("J17").Value = 0.4774346793349168
$ws.Range("K17").Value = 0.09263657957244656
$ws.Range("M17").Value = 0.01187648456057007
$ws.Range("O17").Value = 0.03800475059382423
$ws.Range("S17").Value = 0.1092636579572447
$ws.Range("F18").Value = 0.01687763713080169
$ws.Range("H18").Value = 0.1476793248945148
$ws.Range("I18").Value = 0.1308016877637131
$ws.Range("J18").Value = 0.4388185654008439
$ws.Range("K18").Value = 0.08016877637130802
$ws.Range("M18").Value = 0.01265822784810127
$ws.Range("O18").Value = 0.07172995780590717
$ws.Range("S18").Value = 0.1012658227848101
$ws.Range("F19").Value = 0.01286173633440514
$ws.Range("H19").Value = 0.2033762057877813
$ws.Range("I19").Value = 0.1197749196141479
$ws.Range("J19").Value = 0.3987138263665595
$ws.Range("K19").Value = 0.08038585209003216
$ws.Range("M19").Value = 0.01607717041800643
$ws.Range("N19").Value = 0.0008038585209003215
$ws.Range("O19").Value = 0.06430868167202572
$ws.Range("S19").Value = 0.1036977491961415
